# "Change card class to use index int"
# The Card Class column (C) listed the suit letter for each card, but the
# letters were mis-aligned by one suit-block relative to the computed
# suit index in column E (=QUOTIENT(A,13)). Re-point each 13-row suit
# block at the correct letter, rotating S->H->D->C->S.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# rows 2-14  (suit index 0): was "S", now "H"
$ws.Range("C2:C14").Value = "H"
# rows 15-27 (suit index 1): was "H", now "D"
$ws.Range("C15:C27").Value = "D"
# rows 28-40 (suit index 2): was "D", now "C"
$ws.Range("C28:C40").Value = "C"
# rows 41-53 (suit index 3): was "C", now "S"
$ws.Range("C41:C53").Value = "S"

# Update the view state left by the editing session: active selection and
# scroll position.
$ws.Activate()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 21 } catch {}
try { $win.ScrollColumn = 1 } catch {}
$ws.Range("C33").Select()
